$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells whose new values look like plain numbers stay as text,
# matching the original inline-string (text) representation in the sheet.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated price / volume values
$ws.Range("D2").Value = '72.080.77'
$ws.Range("E2").Value = '  +3.91%  '
$ws.Range("D3").Value = '4.040.63'
$ws.Range("E3").Value = '  +3.65%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '519.41'
$ws.Range("E5").Value = '  -1.91%  '
$ws.Range("D6").Value = '147.31'
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("D7").Value = '0.726'
$ws.Range("E7").Value = '  +18.67%  '
$ws.Range("D8").Value = '4.031.60'
$ws.Range("E8").Value = '  +3.59%  '
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '0.773'
$ws.Range("E10").Value = '  +7.52%  '
$ws.Range("D11").Value = '0.174'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '0.0000326'
$ws.Range("E12").Value = '  -3.14%  '
$ws.Range("D13").Value = '47.60'
$ws.Range("E13").Value = '  +12.81%  '
$ws.Range("D14").Value = '11.07'
$ws.Range("E14").Value = '  +7.72%  '
$ws.Range("D15").Value = '4.689.69'
$ws.Range("E15").Value = '  +3.68%  '
$ws.Range("D16").Value = '4.046.25'
$ws.Range("E16").Value = '  +3.07%  '
$ws.Range("D17").Value = '21.17'
$ws.Range("E17").Value = '  +6.76%  '
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Value = '72.086.67'
$ws.Range("E21").Value = '  +4.05%  '
$ws.Range("D22").Value = '443.06'
$ws.Range("E22").Value = '  +4.18%  '
$ws.Range("D23").Value = '104.84'
$ws.Range("E23").Value = '  +18.91%  '
$ws.Range("E24").Value = '  +4.66%  '
$ws.Range("D25").Value = '14.85'
$ws.Range("E25").Value = '  +4.81%  '
$ws.Range("D26").Value = '3.99'
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("D27").Value = '11.43'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '11.00'
$ws.Range("E28").Value = '  +3.84%  '
$ws.Range("D29").Value = '37.69'
$ws.Range("E29").Value = '  +3.64%  '
$ws.Range("D30").Value = '5.79'
$ws.Range("E30").Value = '  +2.07%  '
$ws.Range("E31").Value = '  +13.99%  '
$ws.Range("D32").Value = '13.65'
$ws.Range("E32").Value = '  +3.43%  '
$ws.Range("E33").Value = '  +3.08%  '
$ws.Range("D34").Value = '679.45'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("E35").Value = '  +15.05%  '
$ws.Range("D36").Value = '66.79'
$ws.Range("E36").Value = '  -3.38%  '
$ws.Range("D37").Value = '42.35'
$ws.Range("E37").Value = '  +5.70%  '
$ws.Range("D38").Value = '0.0₃0862'
$ws.Range("E38").Value = '  -1.87%  '
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("D40").Value = '3.52'
$ws.Range("E40").Value = '  +8.29%  '
$ws.Range("D41").Value = '0.150'
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '0.0498'
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  +1.60%  '
$ws.Range("D46").Value = '0.159'
$ws.Range("E46").Value = '  +12.99%  '
$ws.Range("E47").Value = '  +4.45%  '
$ws.Range("D49").Value = '3.05'
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("D50").Value = '9.14'
$ws.Range("E50").Value = '  +7.45%  '
$ws.Range("E51").Value = '  +2.09%  '
